$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Tshr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 203.7816646666667
$ws.Range("H2").Value = 611.344994
$ws.Range("I2").Value = 0.6667327591988204
$ws.Range("J2").Value = 0.6667327591988205
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.960636
$ws.Range("N2").Value = 2.881908
$ws.Range("O2").Value = 0.2124324572954377
$ws.Range("P2").Value = 0.2124324572954377
$ws.Range("Q2").Value = 195.760003218728
$ws.Range("R2").Value = 1761.840028968552
$ws.Range("S2").Value = 0.1416356783959727
$ws.Range("T2").Value = 0.1416356783959728

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Tshr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 203.7816646666667
$ws.Range("H3").Value = 611.344994
$ws.Range("I3").Value = 0.6667327591988204
$ws.Range("J3").Value = 0.6667327591988205
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.5285266666666667
$ws.Range("N3").Value = 1.58558
$ws.Range("O3").Value = 0.116876963330717
$ws.Range("P3").Value = 0.116876963330717
$ws.Range("Q3").Value = 107.7040439540578
$ws.Range("R3").Value = 969.3363955865201
$ws.Range("S3").Value = 0.07792570024826832
$ws.Range("T3").Value = 0.07792570024826834

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Tshr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 203.7816646666667
$ws.Range("H4").Value = 611.344994
$ws.Range("I4").Value = 0.6667327591988204
$ws.Range("J4").Value = 0.6667327591988205
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.032914666666667
$ws.Range("N4").Value = 9.098744
$ws.Range("O4").Value = 0.6706905793738452
$ws.Range("P4").Value = 0.6706905793738454
$ws.Range("Q4").Value = 618.0523995652818
$ws.Range("R4").Value = 5562.471596087536
$ws.Range("S4").Value = 0.4471713805545793
$ws.Range("T4").Value = 0.4471713805545794

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Tshr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("H5").Value = 189.421768
$ws.Range("I5").Value = 0.2065833519051582
$ws.Range("J5").Value = 0.2065833519051582
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.960636
$ws.Range("N5").Value = 2.881908
$ws.Range("O5").Value = 0.2124324572954377
$ws.Range("P5").Value = 0.2124324572954377
$ws.Range("Q5").Value = 60.655123174816
$ws.Range("R5").Value = 545.896108573344
$ws.Range("S5").Value = 0.04388500908154088
$ws.Range("T5").Value = 0.04388500908154089

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Tshr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 63.14058933333333
$ws.Range("H6").Value = 189.421768
$ws.Range("I6").Value = 0.2065833519051582
$ws.Range("J6").Value = 0.2065833519051582
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.5285266666666667
$ws.Range("N6").Value = 1.58558
$ws.Range("O6").Value = 0.116876963330717
$ws.Range("P6").Value = 0.116876963330717
$ws.Range("Q6").Value = 33.37148521171556
$ws.Range("R6").Value = 300.34336690544
$ws.Range("S6").Value = 0.02414483484535579
$ws.Range("T6").Value = 0.02414483484535579

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Tshr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 63.14058933333333
$ws.Range("H7").Value = 189.421768
$ws.Range("I7").Value = 0.2065833519051582
$ws.Range("J7").Value = 0.2065833519051582
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.032914666666667
$ws.Range("N7").Value = 9.098744
$ws.Range("O7").Value = 0.6706905793738452
$ws.Range("P7").Value = 0.6706905793738454
$ws.Range("Q7").Value = 191.5000194510436
$ws.Range("R7").Value = 1723.500175059392
$ws.Range("S7").Value = 0.1385535079782615
$ws.Range("T7").Value = 0.1385535079782615

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Tshr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 38.719942
$ws.Range("H8").Value = 116.159826
$ws.Range("I8").Value = 0.1266838888960214
$ws.Range("J8").Value = 0.1266838888960214
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.960636
$ws.Range("N8").Value = 2.881908
$ws.Range("O8").Value = 0.2124324572954377
$ws.Range("P8").Value = 0.2124324572954377
$ws.Range("Q8").Value = 37.195770203112
$ws.Range("R8").Value = 334.761931828008
$ws.Range("S8").Value = 0.02691176981792403
$ws.Range("T8").Value = 0.02691176981792404

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Tshr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 38.719942
$ws.Range("H9").Value = 116.159826
$ws.Range("I9").Value = 0.1266838888960214
$ws.Range("J9").Value = 0.1266838888960214
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.5285266666666667
$ws.Range("N9").Value = 1.58558
$ws.Range("O9").Value = 0.116876963330717
$ws.Range("P9").Value = 0.116876963330717
$ws.Range("Q9").Value = 20.46452187878667
$ws.Range("R9").Value = 184.18069690908
$ws.Range("S9").Value = 0.01480642823709292
$ws.Range("T9").Value = 0.01480642823709293

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Gnai2"
$ws.Range("C10").Value = "Tshr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 38.719942
$ws.Range("H10").Value = 116.159826
$ws.Range("I10").Value = 0.1266838888960214
$ws.Range("J10").Value = 0.1266838888960214
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.032914666666667
$ws.Range("N10").Value = 9.098744
$ws.Range("O10").Value = 0.6706905793738452
$ws.Range("P10").Value = 0.6706905793738454
$ws.Range("Q10").Value = 117.4342799842827
$ws.Range("R10").Value = 1056.908519858544
$ws.Range("S10").Value = 0.08496569084100443
$ws.Range("T10").Value = 0.08496569084100446
